$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.636.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.32%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.630.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'212.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.46%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.0624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.84%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.859.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'1.629.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.47%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.18%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'26.632.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.39%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'209.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.82%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E23").Value = "'  +2.70%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +2.76%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.55%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.96%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +4.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.40%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.50%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.166.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.46%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.17%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.806"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.51%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.90%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.770.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.55%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'91.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.70%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'54.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.08%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Mantle"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.409"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.46%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.18%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.04%  "
$ws.Range("E51").Style = "Normal"
